$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---- Sheet1 (List1): rewrite the ranking-overview table (rows 2-19) ----
$ws1.Range("B2").Value = "M21"
$ws1.Range("C2").Value = "M21"
$ws1.Range("D2").Value = "3.4–70–27"
$ws1.Range("E2").Value = "31–33–34–35–36–37–38–39–40–41–42–43–44–45–46–47–48–49–50–51–52–53–54–55–57–59–100––"
$ws1.Range("B3").Value = "M55"
$ws1.Range("C3").Value = "M55 Ž45"
$ws1.Range("D3").Value = "2.6–70–17"
$ws1.Range("E3").Value = "32–49–48–46–42–39–38–35–36–33–50–51–54–55–58–59–100––––––––––––"
$ws1.Range("B4").Value = "W45"
$ws1.Range("C4").Value = "M55 Ž45"
$ws1.Range("D4").Value = "2.6–70–17"
$ws1.Range("E4").Value = "32–49–48–46–42–39–38–35–36–33–50–51–54–55–58–59–100––––––––––––"
$ws1.Range("B5").Value = "W20"
$ws1.Range("C5").Value = "Ž20 Ž35"
$ws1.Range("D5").Value = "2.9–50–17"
$ws1.Range("E5").Value = "33–37–34–38–39–42–41–43–46–47–49–32–51–52–54–57–100––––––––––––"
$ws1.Range("B6").Value = "W35"
$ws1.Range("C6").Value = "Ž20 Ž35"
$ws1.Range("D6").Value = "2.9–50–17"
$ws1.Range("E6").Value = "33–37–34–38–39–42–41–43–46–47–49–32–51–52–54–57–100––––––––––––"
$ws1.Range("B7").Value = "M65"
$ws1.Range("C7").Value = "M65 Ž65"
$ws1.Range("D7").Value = "1.9–50–12"
$ws1.Range("E7").Value = "33–37–36–61–32–51–49–48–53–54–57–100–––––––––––––––––"
$ws1.Range("B8").Value = "W65"
$ws1.Range("C8").Value = "M65 Ž65"
$ws1.Range("D8").Value = "1.9–50–12"
$ws1.Range("E8").Value = "33–37–36–61–32–51–49–48–53–54–57–100–––––––––––––––––"
$ws1.Range("B9").Value = "W16"
$ws1.Range("C9").Value = "Ž16 Ž55"
$ws1.Range("D9").Value = "2–50–13"
$ws1.Range("E9").Value = "34–35–37–48–47–50–51–52–53–54–55–58–100––––––––––––––––"
$ws1.Range("B10").Value = "W55"
$ws1.Range("C10").Value = "Ž16 Ž55"
$ws1.Range("D10").Value = "2–50–13"
$ws1.Range("E10").Value = "34–35–37–48–47–50–51–52–53–54–55–58–100––––––––––––––––"
$ws1.Range("B11").Value = "M20"
$ws1.Range("C11").Value = "Ž21 M35 M20"
$ws1.Range("D11").Value = "3–70–21"
$ws1.Range("E11").Value = "36–35–34–37–49–48–46–45–41–42–40–43–31–32–50–51–53–54–57–59–100––––––––"
$ws1.Range("B12").Value = "M35"
$ws1.Range("C12").Value = "Ž21 M35 M20"
$ws1.Range("D12").Value = "3–70–21"
$ws1.Range("E12").Value = "36–35–34–37–49–48–46–45–41–42–40–43–31–32–50–51–53–54–57–59–100––––––––"
$ws1.Range("B13").Value = "W21"
$ws1.Range("C13").Value = "Ž21 M35 M20"
$ws1.Range("D13").Value = "3–70–21"
$ws1.Range("E13").Value = "36–35–34–37–49–48–46–45–41–42–40–43–31–32–50–51–53–54–57–59–100––––––––"
$ws1.Range("B14").Value = "M16"
$ws1.Range("C14").Value = "M16 M45"
$ws1.Range("D14").Value = "2.8–50–21"
$ws1.Range("E14").Value = "37–36–35–38–39–43–44–45–46–47–49–33–61–32–50–51–53–54–57–58–100––––––––"
$ws1.Range("B15").Value = "M45"
$ws1.Range("C15").Value = "M16 M45"
$ws1.Range("D15").Value = "2.8–50–21"
$ws1.Range("E15").Value = "37–36–35–38–39–43–44–45–46–47–49–33–61–32–50–51–53–54–57–58–100––––––––"
$ws1.Range("B16").Value = "M16"
$ws1.Range("C16").Value = "M16 M45"
$ws1.Range("D16").Value = "2.8–50–21"
$ws1.Range("E16").Value = "37–36–35–38–39–43–44–45–46––49–33–61–32–50–51–53–54–57–58–100––––––––"
$ws1.Range("B17").Value = "M12"
$ws1.Range("C17").Value = "M12 Ž12 OPEN"
$ws1.Range("D17").Value = "1.5–30–12"
$ws1.Range("E17").Value = "61–35–38–49–33–50–51–53–55–58–59–100–––––––––––––––––"
$ws1.Range("B18").Value = "OPEN"
$ws1.Range("C18").Value = "M12 Ž12 OPEN"
$ws1.Range("D18").Value = "1.5–30–12"
$ws1.Range("E18").Value = "61–35–38–49–33–50–51–53–55–58–59–100–––––––––––––––––"
$ws1.Range("B19").Value = "W12"
$ws1.Range("C19").Value = "M12 Ž12 OPEN"
$ws1.Range("D19").Value = "1.5–30–12"
$ws1.Range("E19").Value = "61–35–38–49–33–50–51–53–55–58–59–100–––––––––––––––––"

# Row 14 previously carried an explicit left-aligned override style left over
# from the old table shape; re-level it to the plain (unstyled) look used by the
# rest of the table by copying the style from an unstyled row in the same columns.
$ws1.Range("B14").Style = $ws1.Range("B2").Style
$ws1.Range("C14").Style = $ws1.Range("C2").Style
$ws1.Range("D14").Style = $ws1.Range("D2").Style
$ws1.Range("E14").Style = $ws1.Range("E2").Style

# The table now only spans down to row 19; blank out the old rows 20-25
# (including their leftover override styles) so they read as empty, like the
# rest of the sheet below the table.
$ws1.Range("B20:E25").Clear()

# Drop the trailing block of now-unused pre-formatted rows 234-252.
$ws1.Rows("234:252").Delete()

# ---- Sheet2 (List2): rewrite the start-list table (rows 2-19) ----
$ws2.Range("B2").Value = "M21"
$ws2.Range("C2").Value = "M21"
$ws2.Range("E2").Value = 1
$ws2.Range("G2").Value = "M21"
$ws2.Range("H2").Value = 1
$ws2.Range("B3").Value = "W21"
$ws2.Range("C3").Value = "Ž21 M35 M20"
$ws2.Range("E3").Value = 2
$ws2.Range("G3").Value = "Ž21"
$ws2.Range("H3").Value = 2
$ws2.Range("B4").Value = "M35"
$ws2.Range("C4").Value = "Ž21 M35 M20"
$ws2.Range("E4").Value = 2
$ws2.Range("G4").Value = "M35"
$ws2.Range("H4").Value = 6
$ws2.Range("B5").Value = "M20"
$ws2.Range("C5").Value = "Ž21 M35 M20"
$ws2.Range("E5").Value = 2
$ws2.Range("G5").Value = "M20"
$ws2.Range("H5").Value = 7
$ws2.Range("B6").Value = "M65"
$ws2.Range("C6").Value = "M65 Ž65"
$ws2.Range("E6").Value = 3
$ws2.Range("G6").Value = "M65"
$ws2.Range("H6").Value = 3
$ws2.Range("B7").Value = "W65"
$ws2.Range("C7").Value = "M65 Ž65"
$ws2.Range("E7").Value = 3
$ws2.Range("G7").Value = "Ž65"
$ws2.Range("H7").Value = 10
$ws2.Range("B8").Value = "M55"
$ws2.Range("C8").Value = "M55 Ž45"
$ws2.Range("E8").Value = 4
$ws2.Range("G8").Value = "M55"
$ws2.Range("H8").Value = 4
$ws2.Range("B9").Value = "W45"
$ws2.Range("C9").Value = "M55 Ž45"
$ws2.Range("E9").Value = 4
$ws2.Range("G9").Value = "Ž45"
$ws2.Range("H9").Value = 12
$ws2.Range("B10").Value = "M45"
$ws2.Range("C10").Value = "M16 M45"
$ws2.Range("E10").Value = 5
$ws2.Range("G10").Value = "M45"
$ws2.Range("H10").Value = 5
$ws2.Range("B11").Value = "M16"
$ws2.Range("C11").Value = "M16 M45"
$ws2.Range("E11").Value = 5
$ws2.Range("G11").Value = "M16"
$ws2.Range("H11").Value = 8
$ws2.Range("B12").Value = "M16"
$ws2.Range("C12").Value = "M16 M45"
$ws2.Range("E12").Value = 5
$ws2.Range("G12").Value = "M16"
$ws2.Range("H12").Value = 8
$ws2.Range("B13").Value = "M12"
$ws2.Range("C13").Value = "M12 Ž12 OPEN"
$ws2.Range("E13").Value = 9
$ws2.Range("G13").Value = "M12"
$ws2.Range("H13").Value = 9
$ws2.Range("B14").Value = "W12"
$ws2.Range("C14").Value = "M12 Ž12 OPEN"
$ws2.Range("E14").Value = 9
$ws2.Range("G14").Value = "Ž12"
$ws2.Range("H14").Value = 16
$ws2.Range("B15").Value = "OPEN"
$ws2.Range("C15").Value = "M12 Ž12 OPEN"
$ws2.Range("E15").Value = 9
$ws2.Range("G15").Value = "Otvorena"
$ws2.Range("H15").Value = 17
$ws2.Range("B16").Value = "W55"
$ws2.Range("C16").Value = "Ž16 Ž55"
$ws2.Range("E16").Value = 11
$ws2.Range("G16").Value = "Ž55"
$ws2.Range("H16").Value = 11
$ws2.Range("B17").Value = "W16"
$ws2.Range("C17").Value = "Ž16 Ž55"
$ws2.Range("E17").Value = 11
$ws2.Range("G17").Value = "Ž16"
$ws2.Range("H17").Value = 15
$ws2.Range("B18").Value = "W35"
$ws2.Range("C18").Value = "Ž20 Ž35"
$ws2.Range("E18").Value = 13
$ws2.Range("G18").Value = "Ž35"
$ws2.Range("H18").Value = 13
$ws2.Range("B19").Value = "W20"
$ws2.Range("C19").Value = "Ž20 Ž35"
$ws2.Range("E19").Value = 13
$ws2.Range("G19").Value = "Ž20"
$ws2.Range("H19").Value = 14

# The start-list table shrank from 25 to 19 rows; remove the extra rows entirely.
$ws2.Rows("20:25").Delete()
